$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Convert L2 and O2 from text to real numbers
$ws.Range("L2").Value = 638057
$ws.Range("O2").Value = 9876543210

# Add new row 3 with faculty record (id=14)
$ws.Range("A3").Value = 14
$ws.Range("B3").Value = "man"
$ws.Range("C3").Value = "Doe"
$ws.Range("D3").Value = "hello@gmail.com"
$ws.Range("E3").Value = "Under Graduate"
$ws.Range("F3").Value = "SOftware dev"
$ws.Range("G3").Value = "aws"
$ws.Range("H3").Value = "Chennai"
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = "https://linkedin.com/in/dharunap"
$ws.Range("K3").Value = "tambaram,chennai"
$ws.Range("L3").NumberFormat = "@"
$ws.Range("L3").Value = "638057"
$ws.Range("M3").Value = "chennai"
$ws.Range("N3").Value = "India"
$ws.Range("O3").NumberFormat = "@"
$ws.Range("O3").Value = "9876543210"
$ws.Range("P3").Value = "https://google.com/"
